# Fix duplicate data and correct reconciliation results
# Rows 17-26 were duplicates of rows 7-16 (T006-T015 added twice).
# Renumber the duplicate block to T016-T025 with corrected PV/Delta values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: T006 -> T016
$ws.Range("A17").Value = "T016"
$ws.Range("B17").Value = 92000
$ws.Range("C17").Value = 0.45

# Row 18: T007 -> T017
$ws.Range("A18").Value = "T017"
$ws.Range("B18").Value = -78000
$ws.Range("C18").Value = -0.75

# Row 19: T008 -> T018
$ws.Range("A19").Value = "T018"
$ws.Range("B19").Value = 105000
$ws.Range("C19").Value = 0.58

# Row 20: T009 -> T019 (PV/Delta unchanged, both 0)
$ws.Range("A20").Value = "T019"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0

# Row 21: T010 -> T020
$ws.Range("A21").Value = "T020"
$ws.Range("B21").Value = 82000
$ws.Range("C21").Value = 0.41

# Row 22: T011 -> T021
$ws.Range("A22").Value = "T021"
$ws.Range("B22").Value = -85000
$ws.Range("C22").Value = -0.88

# Row 23: T012 -> T022
$ws.Range("A23").Value = "T022"
$ws.Range("B23").Value = 98000
$ws.Range("C23").Value = 0.52

# Row 24: T013 -> T023 (PV/Delta unchanged, both 0)
$ws.Range("A24").Value = "T023"
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 0

# Row 25: T014 -> T024
$ws.Range("A25").Value = "T024"
$ws.Range("B25").Value = 72000
$ws.Range("C25").Value = 0.39

# Row 26: T015 -> T025
$ws.Range("A26").Value = "T025"
$ws.Range("B26").Value = -92000
$ws.Range("C26").Value = -0.85
